$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''265.42'
$ws.Range("G2").Value = '''9'
$ws.Range("D3").Value = '''22.62'
$ws.Range("G3").Value = '''9'
$ws.Range("D4").Value = '''6.191'
$ws.Range("G4").Value = '''9'
$ws.Range("D5").Value = '''0.06140'
$ws.Range("G5").Value = '''9'
$ws.Range("D6").Value = '''3.563'
$ws.Range("G6").Value = '''9'
$ws.Range("D7").Value = '''6.709'
$ws.Range("G7").Value = '''9'
$ws.Range("D8").Value = '''1.358'
$ws.Range("G8").Value = '''9'
$ws.Range("D9").Value = '''0.8261'
$ws.Range("G9").Value = '''9'
$ws.Range("G10").Value = '''9'
$ws.Range("D11").Value = '''0.1596'
$ws.Range("G11").Value = '''9'
$ws.Range("D12").Value = '''0.08244'
$ws.Range("G12").Value = '''9'
$ws.Range("D13").Value = '''0.03401'
$ws.Range("G13").Value = '''9'
$ws.Range("D14").Value = '''0.03163'
$ws.Range("G14").Value = '''9'
$ws.Range("B15").Value = 'ProBitToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D15").Value = '''0.1208'
$ws.Range("E15").Value = '14ProBitTokenPROB'
$ws.Range("G15").Value = '''9'
$ws.Range("B16").Value = 'BitMartToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D16").Value = '''0.09238'
$ws.Range("E16").Value = '15BitMartTokenBMX'
$ws.Range("G16").Value = '''9'
$ws.Range("B17").Value = 'MCDex'
$ws.Range("C17").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D17").Value = '''3.895'
$ws.Range("E17").Value = '16MCDexMCB'
$ws.Range("G17").Value = '''9'
$ws.Range("B18").Value = 'BitForexToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D18").Value = '''0.001688'
$ws.Range("E18").Value = '17BitForexTokenBF'
$ws.Range("G18").Value = '''9'
$ws.Range("B19").Value = 'CoinExToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D19").Value = '''0.04811'
$ws.Range("E19").Value = '18CoinExTokenCET'
$ws.Range("G19").Value = '''9'
$ws.Range("B20").Value = 'TigerCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D20").Value = '''0.006231'
$ws.Range("E20").Value = '19TigerCashTCH'
$ws.Range("G20").Value = '''9'
$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D21").Value = '''0.006282'
$ws.Range("E21").Value = '20HotbitTokenHTB'
$ws.Range("G21").Value = '''9'
$ws.Range("B22").Value = 'BitKan'
$ws.Range("C22").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D22").Value = '''0.001097'
$ws.Range("E22").Value = '21BitKanKAN'
$ws.Range("G22").Value = '''9'
$ws.Range("B23").Value = 'NitroEx'
$ws.Range("C23").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D23").Value = '''0.0001501'
$ws.Range("E23").Value = '22NitroExNTX'
$ws.Range("G23").Value = '''9'
$ws.Range("B24").Value = 'LEO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D24").Value = '''3.718'
$ws.Range("E24").Value = '23LEOLEO'
$ws.Range("G24").Value = '''9'
$ws.Range("B25").Value = 'BTSEToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D25").Value = '''2.301'
$ws.Range("E25").Value = '24BTSETokenBTSE'
$ws.Range("G25").Value = '''9'
$ws.Range("B26").Value = 'BitpandaEcosystemToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D26").Value = '''0.3382'
$ws.Range("E26").Value = '25BitpandaEcosystemTokenBEST'
$ws.Range("G26").Value = '''9'
$ws.Range("D27").Value = '''0.0002684'
$ws.Range("G27").Value = '''9'
$ws.Range("G28").Value = '''9'
$ws.Range("G29").Value = '''9'
$ws.Range("G30").Value = '''9'
$ws.Range("G31").Value = '''9'
$ws.Range("G32").Value = '''9'
$ws.Range("G33").Value = '''9'
$ws.Range("G34").Value = '''9'
$ws.Range("G35").Value = '''9'
$ws.Range("G36").Value = '''9'
$ws.Range("G37").Value = '''9'
$ws.Range("G38").Value = '''9'
$ws.Range("G39").Value = '''9'
$ws.Range("D40").Value = '''0.04605'
$ws.Range("G40").Value = '''9'
$ws.Range("D41").Value = '''0.006970'
$ws.Range("G41").Value = '''9'
$ws.Range("D42").Value = '''0.1131'
$ws.Range("G42").Value = '''9'
$ws.Range("D43").Value = '''0.003248'
$ws.Range("G43").Value = '''9'
$ws.Range("D44").Value = '''0.01093'
$ws.Range("G44").Value = '''9'
$ws.Range("D45").Value = '''0.00006163'
$ws.Range("G45").Value = '''9'
$ws.Range("G46").Value = '''9'
$ws.Range("D47").Value = '''0.7706'
$ws.Range("G47").Value = '''9'
$ws.Range("D48").Value = '''0.2051'
$ws.Range("G48").Value = '''9'
$ws.Range("D49").Value = '''0.00002102'
$ws.Range("G49").Value = '''9'
$ws.Range("D50").Value = '''0.01241'
$ws.Range("G50").Value = '''9'
$ws.Range("G51").Value = '''9'

Write-Host "Applied cryptos.xlsx update"
